# Apply the ValueSet-fr-mp-dose-form.xlsx changes described in the commit diff.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1 ("Metadata"): update the URL, Date, Jurisdiction and Description
# value cells (column B) in place. Row/column layout is unchanged.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Row 2: URL
$ws1.Range("B2").Value = "https://hl7.fr/ig/fhir/medication/ValueSet/fr-mp-dose-form"

# Row 8: Date
$ws1.Range("B8").Value = "2026-01-15T08:54:26+00:00"

# Row 11: Jurisdiction
$ws1.Range("B11").Value = "FRANCE"

# Row 12: Description
$ws1.Range("B12").Value = "Le jeu de valeurs à utiliser pour coder l'élément *doseForm* des ressources *FrMedication*."

# ---------------------------------------------------------------------------
# Sheet 2 ("Include #0"): the constraint/Operation table is replaced with a
# "Codes / All codes" summary, and the System URI value changes. Column C is
# no longer used, so it is fully cleared (not just its contents) to shrink
# the sheet dimension back down to A1:B4.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Row 1: Property/Operation/Value header -> just "Codes"
$ws2.Range("A1").Value = "Codes"
$ws2.Range("B1:C1").Clear()

# Row 2: constraint / = / (expression) -> just "All codes"
$ws2.Range("A2").Value = "All codes"
$ws2.Range("B2:C2").Clear()

# Row 3 stays the blank separator row (A3/B3 already empty) - leave as is,
# but make sure any leftover column C content is gone.
$ws2.Range("C3").Clear()

# Row 4: System URI value changes; label stays the same.
$ws2.Range("A4").Value = "System URI"
$ws2.Range("B4").Value = "http://standardterms.edqm.eu"
$ws2.Range("C4").Clear()
